$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (SRE subject) - updated Note line
$ws.Range("D2").Value = 'Subject: Software Requirement Engineering;
Instructor: Mr Fazal Wahab;
ClassSenior: Fakhar | +92 ‭332 0278998‬;
Note: From this Saturday (30 Nov 2019) to end of semester means four weeks, there will extra 1 and half hour makeup class of SRE after maghrib prayer.;
CreditHours: 3.0;'

# Row 11 D - new SRE Week 8 content (previously empty, needs style to match other data cells)
$ws.Range("D11").Value = 'Topic: name- Week 8, lectures- 3 Lectures, duration- 01:22;
Video: link- https://drive.google.com/file/d/18sP4Wf0WZh_3gh0PCDZaPH0YT7HsDmhn/preview, name- SRE Week8 Part#1.mp4, duration- 00:01;
Video: link- https://drive.google.com/file/d/1cFUxdRC5VwlsjxILvMbMqO6SVWkXwZB1/preview, name- SRE Week8 Part#2.mp4, duration- 00:51;
Video: link- https://drive.google.com/file/d/128EIgEoYNkqrN_PCG30nwpn1MGaNhK-J/preview, name- SRE Week8 Part#3.mp4, duration- 00:30;
Assignment: name- Assignment no 3 (Deadline 21 Dec 2019), img- /SRE/Assignment 3.png;'
$ws.Range("D11").WrapText = $true
$ws.Range("D11").Font.Bold = $true

# Row 11 E - ALGO Week 8 (assignment moved after videos)
$ws.Range("E11").Value = 'Topic: name- Week 8, lectures- 2 Lectures, duration- 01:33;
Video: link- https://drive.google.com/file/d/1lpPNJAvs6WzQuJ6z0dhX96HCxnarShVs/preview, name- AD&AA Week # 8 Part1.mp4, duration- 00:54;
Video: link- https://drive.google.com/file/d/1A8aBmhTcErBz6hETG1YYLN-GMCLMZUSX/preview, name- AD&AA Week # 8 Part2.mp4, duration- 00:39;
Assignment: name- Assignment 2 (Deadline 22 Nov), img- /ALGO/Assignment 2.png*/ALGO/Assignment 2 DL.png;
Slides: slide- lec4.pptx;
Topics: Redex Sort, Bubble Sort;'

# Row 12 E - ALGO Week 9 (videos added)
$ws.Range("E12").Value = 'Topic: name- Week 9, lectures- 2 Lectures, duration- 01:51;
Video: link- https://drive.google.com/file/d/1430C-n3l2kRoy2Qn76kBep_umfClXmkF/preview, name- AD&AA Week # 9 Part1.mp4, duration- 00:38;
Video: link- https://drive.google.com/file/d/1UBidBALGcA7KB8t11RDEwlTti64BilJf/preview, name- AD&AA Week # 9 Part2.mp4, duration- 01:13;
Assignment: name- Assignment 3 (Deadline 29 Nov), img- /ALGO/Assignment 3.png;
Slides: slide- lec5.ppt;
Important: Quiz in next class from lecture 5 (above slides);'

# Row 10 C - TPL Week 7 (Assignment Only wording)
$ws.Range("C10").Value = 'Topic: name- Week 7, lectures- Assignment Only, duration- ;
Assignment: name- Assignment no 2 (Deadline 23 Nov), img- /TPL/Assignment2.jpeg*/TPL/mid20191.jpg*/TPL/mid20192.jpg;
Assignment: name- Assignment no 2 Solved, link- Assignment 2 Solved.docx;
Note: heading- Note, text- No lecture conducted during this week due to exam on 17 Nov 2019. Sir only gave above assignment.;'

# Row 10 D - SRE Week 7 (Assignment Only wording)
$ws.Range("D10").Value = 'Topic: name- Week 7, lectures- Assignment Only, duration- ;
Assignment: name- Assignment no 2 (Deadline 23 Nov), link- Assignment2.pdf;
Assignment: name- Assignment no 2 Solved, link- Assignment 2 Solved.docx;
Books: name- Engineering and Managing software requirements by Aybuke Aurum and Claes Wohlin, link- Engineering and Managing software requirements.pdf;
Note: heading- Note, text- No lecture conducted during this week due to exam on 17 Nov 2019. Sir only gave above assignment;'

# Row 11 C - TPL Week 8 (duration blank)
$ws.Range("C11").Value = 'Topic: name- Week 8, lectures- , duration- ;
Assignment: name- Assignment no 3 (Deadline 15 Dec 2019), img- /TPL/Assignment No 3 TPL.png;
Note: heading- Note, text- Video lectures not uploaded on google drive yet. If you want them faster please call Kamran (Abasyn) at 0302 5003156‬.;'

# Row heights
$ws.Rows.Item(2).RowHeight = 119
$ws.Rows.Item(9).RowHeight = 409.5
$ws.Rows.Item(12).RowHeight = 238

# View / selection
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("D10").Select()
